$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 ("Lens"): replace the old plain-text URL with a real hyperlink
# pointing at the new Alibaba lens-holder-stock listing.
$lensUrl = "https://www.alibaba.com/product-detail/1mm-1-5mm-2mm-4mm-12mm_60275500727.html?spm=a2700.8443308.0.0.4b7a3e5fd82dtX"
$ws.Hyperlinks.Add($ws.Range("D8"), $lensUrl)
$ws.Range("D8").Value = $lensUrl

# Row 10: new BOM line "Lens Holder", quantity 1
$ws.Range("B10").Value = "Lens Holder"
$ws.Range("C10").Value = 1

# Move the active selection to D10, matching the saved cursor position.
$ws.Range("D10").Select()
